# Add a new "J.Fromm" sheet (Week 16 simulated row), matching the layout
# already used by the "M.Glennon" and "D.Jones" sheets.

$wb = $excel.ActiveWorkbook

# Use the existing "M.Glennon" sheet as the formatting/layout template.
$template = $wb.Worksheets.Item("M.Glennon")

# New sheet goes after the last existing sheet ("D.Jones").
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "J.Fromm"

# Column headers (row 1, columns B:G).
$headers = @("Short Att", "Short Comp", "Deep Att", "Deep Comp", "Short Int", "Deep Int")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, 2 + $i).Value = $headers[$i]
}

# Row labels (Home / Road).
$ws.Cells.Item(2, 1).Value = "H"
$ws.Cells.Item(3, 1).Value = "R"

# Simulated Week 16 stats - all zeros (placeholder, not yet logged).
for ($r = 2; $r -le 3; $r++) {
    for ($c = 2; $c -le 7; $c++) {
        $ws.Cells.Item($r, $c).Value = 0
    }
}

# Copy the header/label formatting (bold, centered, bordered) from the
# template sheet, without touching A1 (which stays empty on every sheet).
$template.Range("B1:G1").Copy()
$ws.Range("B1:G1").PasteSpecial(-4122)
$template.Range("A2:A3").Copy()
$ws.Range("A2:A3").PasteSpecial(-4122)
